$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.637.02"
$ws.Range("E2").Value = "  -1.95%  "
$ws.Range("D3").Value = "2.418.84"
$ws.Range("E3").Value = "  +5.57%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.49"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.91"
$ws.Range("E6").Value = "  -3.31%  "
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -0.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.52"
$ws.Range("E10").Value = "  -5.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0792"
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.18"
$ws.Range("E12").Value = "  -1.71%  "
$ws.Range("E13").Value = "  +0.73%  "
$ws.Range("D14").Value = "2.789.63"
$ws.Range("E14").Value = "  +5.67%  "
$ws.Range("D15").Value = "2.431.93"
$ws.Range("E15").Value = "  +6.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.34"
$ws.Range("E16").Value = "  +4.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.842"
$ws.Range("E17").Value = "  +4.48%  "
$ws.Range("D18").Value = "45.629.68"
$ws.Range("E18").Value = "  -1.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.11"
$ws.Range("E19").Value = "  +1.84%  "
$ws.Range("D20").Value = "0.0₃0952"
$ws.Range("E20").Value = "  +2.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.17"
$ws.Range("E21").Value = "  +2.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.63"
$ws.Range("E22").Value = "  +2.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "244.22"
$ws.Range("E23").Value = "  -1.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.80"
$ws.Range("E24").Value = "  -2.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.94"
$ws.Range("E25").Value = "  +1.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "38.94"
$ws.Range("E27").Value = "  -8.33%  "
$ws.Range("E28").Value = "  -1.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.82"
$ws.Range("E29").Value = "  +0.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.87"
$ws.Range("E30").Value = "  +20.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "21.40"
$ws.Range("E31").Value = "  +7.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.59"
$ws.Range("E32").Value = "  -0.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.74"
$ws.Range("E33").Value = "  -2.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "148.15"
$ws.Range("E34").Value = "  +0.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0778"
$ws.Range("E35").Value = "  -1.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.99"
$ws.Range("E36").Value = "  +12.94%  "
$ws.Range("E37").Value = "  -0.88%  "
$ws.Range("E38").Value = "  -0.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.23"
$ws.Range("E39").Value = "  -3.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.88"
$ws.Range("E40").Value = "  -2.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0301"
$ws.Range("E41").Value = "  +0.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.28"
$ws.Range("E42").Value = "  -1.44%  "
$ws.Range("D43").Value = "1.962.88"
$ws.Range("E43").Value = "  +7.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.82"
$ws.Range("E45").Value = "  +4.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.76"
$ws.Range("E46").Value = "  -9.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.70"
$ws.Range("E47").Value = "  +11.11%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "100.38"
$ws.Range("E48").Value = "  +5.37%  "
$ws.Range("D49").Value = "2.656.95"
$ws.Range("E49").Value = "  +5.61%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.186"
$ws.Range("E50").Value = "  -3.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "14.85"
$ws.Range("E51").Value = "  +10.95%  "
